# Update countries & provincias Spain
# Refresh COVID country stats that moved a handful of countries past their
# neighbours in the (descending, by "Casos totales") ranking, plus a few
# standalone numeric refreshes and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-Row($Row, $Name, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value = $Name
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# --- Standalone numeric refreshes (country identity unchanged) ---
Set-Row 21  "Pakistan"      306304 633 292869 7015  0 4  6420
Set-Row 35  "Kazajistan"    107307 45  101941 3695  0 0  1671
Set-Row 50  "Honduras"      71616  473 22065  47367 0 18 2184
Set-Row 78  "Australia"     26912  14  24063  1998  0 2  851
Set-Row 156 "Nueva Zelanda" 1815   0   1728   62    0 0  25
Set-Row 158 "Belice"        1627   21  918    688   0 1  21

# --- Belgica overtakes Egipto and Marruecos in the ranking ---
Set-Row 37 "Belgica"   102295 1547 18965 73382 0 4 9948
Set-Row 38 "Egipto"    102015 0    89532 6713  0 0 5770
Set-Row 39 "Marruecos" 101743 0    80732 19181 0 0 1830

# --- Timor Oriental overtakes Santa Lucia (tied totals, simple swap) ---
Set-Row 204 "Timor Oriental" 27 0 26 1 0 0 0
Set-Row 205 "Santa Lucia"    27 0 26 1 0 0 0

# --- Islas Malvinas overtakes Montserrat (tied totals, simple swap) ---
Set-Row 214 "Islas Malvinas" 13 0 13 0 0 0 0
Set-Row 215 "Montserrat"     13 0 12 0 0 0 1

# --- Update the "last refreshed" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 04:57"
